# TrialsSetup.xlsx update
# - Rename the "Trial Name (GIT)" column header to "Trial Name"
#   (both the worksheet cell and the backing table's column name)
# - Update BNT323-01's Progress value from 88 to 100
# - Update the active selection to E3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell A1 (this also renames the backing table's column
# header, since A1 is the header cell of the Query1 ListObject)
$ws.Range("A1").Value = "Trial Name"

# Update BNT323-01's progress value (row 6, column B)
$ws.Range("B6").Value = 100

# Move the active selection to E3
$ws.Range("E3").Select()
